# Applies the July 28, 2020 10:32:58 PM PT run results to the covid
# disparities output workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Texas -- Bexar County) ---
$ws.Range("B3").Value = 44040
$ws.Range("C3").Value = 37984
$ws.Range("D3").Value = 335

# --- Row 6 (California - San Francisco) -- error message changed ---
$ws.Range("O6").Value = "An error occurred. ... KeyError('Date_Uploaded.Data as of')"

# --- Row 7 (Tennessee) ---
$ws.Range("B7").Value = 44040
$ws.Range("C7").Value = 99044
$ws.Range("D7").Value = 999
$ws.Range("E7").Value = 18749
$ws.Range("F7").Value = 343
$ws.Range("G7").Value = 18.93
$ws.Range("H7").Value = 34.33

# --- Row 11 (California - San Diego) ---
$ws.Range("B11").Value = 44040
$ws.Range("C11").Value = 28005
$ws.Range("D11").Value = 547
$ws.Range("E11").Value = 1034
$ws.Range("G11").Value = 4.73
$ws.Range("H11").Value = 3.77
$ws.Range("K11").Value = 21845
$ws.Range("L11").Value = 531

# --- Row 36 (Washington) ---
$ws.Range("B36").Value = 44040
$ws.Range("C36").Value = 54205
$ws.Range("D36").Value = 1548
$ws.Range("E36").Value = 1984
$ws.Range("G36").Value = 5.52
$ws.Range("H36").Value = 3.26
$ws.Range("K36").Value = 35958
$ws.Range("L36").Value = 1505

# --- Row 39 (Delaware) -- error message changed ---
$ws.Range("O39").Value = "An error occurred. ... WebDriverException('unknown error: session deleted because of page crash`nfrom unknown error: cannot determine loading status`nfrom tab crashed`n  (Session info: headless chrome=83.0.4103.116)', None, None)"

# --- Row 41 (Iowa) ---
$ws.Range("B41").Value = 44041
$ws.Range("C41").Value = 42928
$ws.Range("E41").Value = 3407
$ws.Range("G41").Value = 7.94

# --- Row 45 (Ohio) -- run failed this time; most fields go blank ---
$ws.Range("B45:H45").Clear()
$ws.Range("J45").Value = $false
$ws.Range("K45:L45").Clear()
$ws.Range("O45").Value = "An error occurred. ... AttributeError(""'NoneType' object has no attribute 'body'"")"

# --- Row 50 (NewYork) -- run failed this time; most fields go blank ---
$ws.Range("B50:D50").Clear()
$ws.Range("F50").Clear()
$ws.Range("H50").Clear()
$ws.Range("L50").Clear()
$ws.Range("O50").Value = "An error occurred. ... ConnectionRefusedError(111, 'Connection refused')"

# --- Insert two new rows (Wyoming, SouthDakota) ahead of the existing
#     SouthCarolina row, which shifts SouthCarolina from row 51 to row 53 ---
$ws.Rows.Item(51).Insert()
$ws.Rows.Item(51).Insert()

# Row 51 -- Wyoming (errored before any data besides population was fetched)
$ws.Range("A51").Value = "Wyoming"
$ws.Range("B51:H51").Clear()
$ws.Range("I51").Value = $false
$ws.Range("J51").Value = $false
$ws.Range("K51:L51").Clear()
$ws.Range("M51").Value = 5540
$ws.Range("N51").Value = 0.95
$ws.Range("O51").Value = "An error occurred. ... JSONDecodeError('Expecting value: line 1 column 1 (char 0)')"

# Row 52 -- SouthDakota
$ws.Range("A52").Value = "SouthDakota"
$ws.Range("B52").Value = 44040
$ws.Range("B52").NumberFormat = "YYYY-MM-DD"
$ws.Range("C52").Value = 8492
$ws.Range("D52").Value = 123
$ws.Range("E52").Value = 1008
$ws.Range("F52").Clear()
$ws.Range("G52").Value = 11.87
$ws.Range("H52").Clear()
$ws.Range("I52").Value = $false
$ws.Range("J52").Value = $false
$ws.Range("K52").Value = 8492
$ws.Range("L52").Clear()
$ws.Range("M52").Clear()
$ws.Range("N52").Clear()
$ws.Range("O52").Value = "Success!"
